$d = $word.ActiveDocument

# Locate the paragraph that ends the "Commit 15" block
# (the one containing the set()/get()/keys()... description).
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*set(), get(), keys()*Object.assign(), optional chaining*") {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq 0) {
    throw "Could not find target paragraph"
}

$targetPara = $d.Paragraphs.Item($targetIndex)

# Insert a new empty paragraph right after it, inheriting the same
# paragraph/run formatting (ListParagraph style, justified, sz 24).
$targetPara.Range.InsertParagraphAfter()

# The newly created paragraph is now the one right after $targetIndex.
$commitPara = $d.Paragraphs.Item($targetIndex + 1)
$commitPara.Range.InsertBefore("Commit 16:")

# Insert another new paragraph (with the same inherited formatting)
# after the "Commit 16:" paragraph, and fill it with the description text.
$commitPara2 = $d.Paragraphs.Item($targetIndex + 1)
$commitPara2.Range.InsertParagraphAfter()
$descPara = $d.Paragraphs.Item($targetIndex + 2)
$descPara.Range.InsertBefore("Methods, this keyword, global object")
